$data = @(
    ,("ECs", "Efna2", "Epha3", "ECs", 3, 1, 1.417859666666667, 4.253579, 0.1472651073415806, 0.1472651073415806, 2, 0.6666666666666666, 0.1220006666666667, 0.366002, 0.003000890193633104, 0.003000890193633103, 0.1729798245731111, 1.556818421158, 0.0004419264164856756, 0.0004419264164856757)
    ,("ECs", "Efna2", "Epha3", "FAPs", 3, 1, 1.417859666666667, 4.253579, 0.1472651073415806, 0.1472651073415806, 3, 1, 40.12734033333333, 120.382021, 0.987025279393649, 0.987025279393649, 56.89493738923989, 512.054436503159, 0.1453543837187593, 0.1453543837187593)
    ,("ECs", "Efna2", "Epha3", "sCs", 3, 1, 1.417859666666667, 4.253579, 0.1472651073415806, 0.1472651073415806, 3, 1, 0.4054843333333333, 1.216453, 0.009973830412717882, 0.00997383041271788, 0.5749198816985556, 5.174278935287, 0.00146879720633562, 0.00146879720633562)
    ,("FAPs", "Efna2", "Epha3", "ECs", 3, 1, 6.189892666666666, 18.569678, 0.6429093297593833, 0.6429093297593833, 2, 0.6666666666666666, 0.1220006666666667, 0.366002, 0.003000890193633104, 0.003000890193633103, 0.7551710319284444, 6.796539287356, 0.001929300303070165, 0.001929300303070165)
    ,("FAPs", "Efna2", "Epha3", "FAPs", 3, 1, 6.189892666666666, 18.569678, 0.6429093297593833, 0.6429093297593833, 3, 1, 40.12734033333333, 120.382021, 0.987025279393649, 0.987025279393649, 248.3839296621375, 2235.455366959238, 0.6345677608305389, 0.6345677608305389)
    ,("FAPs", "Efna2", "Epha3", "sCs", 3, 1, 6.189892666666666, 18.569678, 0.6429093297593833, 0.6429093297593833, 3, 1, 0.4054843333333333, 1.216453, 0.009973830412717882, 0.00997383041271788, 2.509904501348222, 22.589140512134, 0.006412268625774207, 0.006412268625774206)
    ,("M2", "Efna2", "Epha3", "ECs", 3, 1, 1.012114666666667, 3.036344, 0.105122656728831, 0.105122656728831, 2, 0.6666666666666666, 0.1220006666666667, 0.366002, 0.003000890193633104, 0.003000890193633103, 0.1234786640764444, 1.111307976688, 0.000315461549706208, 0.0003154615497062079)
    ,("M2", "Efna2", "Epha3", "FAPs", 3, 1, 1.012114666666667, 3.036344, 0.105122656728831, 0.105122656728831, 3, 1, 40.12734033333333, 120.382021, 0.987025279393649, 0.987025279393649, 40.61346968569156, 365.521227171224, 0.1037587196283771, 0.1037587196283771)
    ,("M2", "Efna2", "Epha3", "sCs", 3, 1, 1.012114666666667, 3.036344, 0.105122656728831, 0.105122656728831, 3, 1, 0.4054843333333333, 1.216453, 0.009973830412717882, 0.00997383041271788, 0.4103966408702223, 3.693569767832, 0.001048475550747717, 0.001048475550747717)
    ,("sCs", "Efna2", "Epha3", "ECs", 3, 1, 1.008073333333333, 3.02422, 0.1047029061702051, 0.1047029061702051, 2, 0.6666666666666666, 0.1220006666666667, 0.366002, 0.003000890193633104, 0.003000890193633103, 0.1229856187155556, 1.10687056844, 0.0003142019243710556, 0.0003142019243710555)
    ,("sCs", "Efna2", "Epha3", "FAPs", 3, 1, 1.008073333333333, 3.02422, 0.1047029061702051, 0.1047029061702051, 3, 1, 40.12734033333333, 120.382021, 0.987025279393649, 0.987025279393649, 40.45130172762444, 364.06171554862, 0.1033444152159737, 0.1033444152159737)
    ,("sCs", "Efna2", "Epha3", "sCs", 3, 1, 1.008073333333333, 3.02422, 0.1047029061702051, 0.1047029061702051, 3, 1, 0.4054843333333333, 1.216453, 0.009973830412717882, 0.00997383041271788, 0.4087579435177778, 3.67882149166, 0.001044289029860339, 0.001044289029860338)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}
